$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111964622
$ws.Range("B2").Value = 89845
$ws.Range("C2").Value = 'Ovaliderad'
$ws.Range("D2").Value = 'VU'
$ws.Range("E2").Value = 1209
$ws.Range("F2").Value = 'Rynkskinn'
$ws.Range("G2").Value = 'Phlebia centrifuga'
$ws.Range("H2").Value = 'P.Karst.'
$ws.Range("P2").Value = 'Kallhögen 5, Vb'
$ws.Range("Q2").Value = 734972.3834676194
$ws.Range("R2").Value = 7088252.533270728
$ws.Range("S2").Value = 20
$ws.Range("T2").Value = 'Västerbotten'
$ws.Range("U2").Value = 'Vännäs'
$ws.Range("V2").Value = 'Västerbotten'
$ws.Range("W2").Value = 'Vännäs'
$ws.Range("Z2").Value = '16:12'
$ws.Range("AB2").Value = '16:12'
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AW2").Value = 'Billy Lindblom'
$ws.Range("AX2").Value = 'Billy Lindblom'

# Row 3
$ws.Range("A3").Value = 111964457
$ws.Range("B3").Value = 56398
$ws.Range("C3").Value = 'Ovaliderad'
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = 'Tretåig hackspett'
$ws.Range("G3").Value = 'Picoides tridactylus'
$ws.Range("H3").Value = '(Linnaeus, 1758)'
$ws.Range("P3").Value = 'Kallhögen 5, Vb'
$ws.Range("Q3").Value = 734949.4564622594
$ws.Range("R3").Value = 7088268.525185317
$ws.Range("S3").Value = 20
$ws.Range("T3").Value = 'Västerbotten'
$ws.Range("U3").Value = 'Vännäs'
$ws.Range("V3").Value = 'Västerbotten'
$ws.Range("W3").Value = 'Vännäs'
$ws.Range("Z3").Value = '16:01'
$ws.Range("AB3").Value = '16:01'
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = 'Billy Lindblom'
$ws.Range("AX3").Value = 'Billy Lindblom'

# Row 4
$ws.Range("A4").Value = 111965439
$ws.Range("B4").Value = 56398
$ws.Range("C4").Value = 'Ovaliderad'
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = 'Tretåig hackspett'
$ws.Range("G4").Value = 'Picoides tridactylus'
$ws.Range("H4").Value = '(Linnaeus, 1758)'
$ws.Range("P4").Value = 'Kallhögen 5, Vb'
$ws.Range("Q4").Value = 734926.7697699566
$ws.Range("R4").Value = 7088234.05367971
$ws.Range("S4").Value = 20
$ws.Range("T4").Value = 'Västerbotten'
$ws.Range("U4").Value = 'Vännäs'
$ws.Range("V4").Value = 'Västerbotten'
$ws.Range("W4").Value = 'Vännäs'
$ws.Range("Z4").Value = '16:40'
$ws.Range("AB4").Value = '16:40'
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = 'Billy Lindblom'
$ws.Range("AX4").Value = 'Billy Lindblom'

# Row 6
$ws.Range("A6").Value = 111964863
$ws.Range("B6").Value = 89745
$ws.Range("C6").Value = 'Ovaliderad'
$ws.Range("D6").Value = 'VU'
$ws.Range("E6").Value = 2062
$ws.Range("F6").Value = 'Ulltickeporing'
$ws.Range("G6").Value = 'Skeletocutis brevispora'
$ws.Range("H6").Value = 'Niemelä'
$ws.Range("P6").Value = 'Kallhögen 5, Vb'
$ws.Range("Q6").Value = 734972.3834676194
$ws.Range("R6").Value = 7088252.533270728
$ws.Range("S6").Value = 20
$ws.Range("T6").Value = 'Västerbotten'
$ws.Range("U6").Value = 'Vännäs'
$ws.Range("V6").Value = 'Västerbotten'
$ws.Range("W6").Value = 'Vännäs'
$ws.Range("Z6").Value = '16:12'
$ws.Range("AB6").Value = '16:12'
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = 'Billy Lindblom'
$ws.Range("AX6").Value = 'Billy Lindblom'

# Row 7
$ws.Range("A7").Value = 111964847
$ws.Range("B7").Value = 89405
$ws.Range("C7").Value = 'Ovaliderad'
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = 'Ullticka'
$ws.Range("G7").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H7").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("P7").Value = 'Kallhögen 5, Vb'
$ws.Range("Q7").Value = 734972.3834676194
$ws.Range("R7").Value = 7088252.533270728
$ws.Range("S7").Value = 20
$ws.Range("T7").Value = 'Västerbotten'
$ws.Range("U7").Value = 'Vännäs'
$ws.Range("V7").Value = 'Västerbotten'
$ws.Range("W7").Value = 'Vännäs'
$ws.Range("Y7").Value = '''2023-09-08'
$ws.Range("Z7").Value = '16:12'
$ws.Range("AA7").Value = '''2023-09-08'
$ws.Range("AB7").Value = '16:12'
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = 'Billy Lindblom'
$ws.Range("AX7").Value = 'Billy Lindblom'

# Row 8
$ws.Range("A8").Value = 111964175
$ws.Range("B8").Value = 89423
$ws.Range("C8").Value = 'Ovaliderad'
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 5432
$ws.Range("F8").Value = 'Granticka'
$ws.Range("G8").Value = 'Porodaedalea chrysoloma'
$ws.Range("H8").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("P8").Value = 'Kallhögen 5, Vb'
$ws.Range("Q8").Value = 734896.4627943118
$ws.Range("R8").Value = 7088342.483217424
$ws.Range("S8").Value = 20
$ws.Range("T8").Value = 'Västerbotten'
$ws.Range("U8").Value = 'Vännäs'
$ws.Range("V8").Value = 'Västerbotten'
$ws.Range("W8").Value = 'Vännäs'
$ws.Range("Y8").Value = '''2023-09-08'
$ws.Range("Z8").Value = '15:42'
$ws.Range("AA8").Value = '''2023-09-08'
$ws.Range("AB8").Value = '15:42'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = 'Billy Lindblom'
$ws.Range("AX8").Value = 'Billy Lindblom'

# Row 9
$ws.Range("A9").Value = 111964050
$ws.Range("B9").Value = 90065
$ws.Range("C9").Value = 'Ovaliderad'
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 898
$ws.Range("F9").Value = 'Blackticka'
$ws.Range("G9").Value = 'Steccherinum collabens'
$ws.Range("H9").Value = '(Fr.) Vesterholt'
$ws.Range("P9").Value = 'Kallhögen 5, Vb'
$ws.Range("Q9").Value = 734893.3330648565
$ws.Range("R9").Value = 7088354.646951701
$ws.Range("S9").Value = 20
$ws.Range("T9").Value = 'Västerbotten'
$ws.Range("U9").Value = 'Vännäs'
$ws.Range("V9").Value = 'Västerbotten'
$ws.Range("W9").Value = 'Vännäs'
$ws.Range("Y9").Value = '''2023-09-08'
$ws.Range("Z9").Value = '15:42'
$ws.Range("AA9").Value = '''2023-09-08'
$ws.Range("AB9").Value = '15:42'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AW9").Value = 'Billy Lindblom'
$ws.Range("AX9").Value = 'Billy Lindblom'

# Row 10
$ws.Range("A10").Value = 111965370
$ws.Range("B10").Value = 81248
$ws.Range("C10").Value = 'Ovaliderad'
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 1312
$ws.Range("F10").Value = 'Gammelgransskål'
$ws.Range("G10").Value = 'Pseudographis pinicola'
$ws.Range("H10").Value = '(Nyl.) Rehm'
$ws.Range("P10").Value = 'Kallhögen 5, Vb'
$ws.Range("Q10").Value = 734939.7547518623
$ws.Range("R10").Value = 7088232.371273324
$ws.Range("S10").Value = 20
$ws.Range("T10").Value = 'Västerbotten'
$ws.Range("U10").Value = 'Vännäs'
$ws.Range("V10").Value = 'Västerbotten'
$ws.Range("W10").Value = 'Vännäs'
$ws.Range("Y10").Value = '''2023-09-08'
$ws.Range("Z10").Value = '16:38'
$ws.Range("AA10").Value = '''2023-09-08'
$ws.Range("AB10").Value = '16:38'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AW10").Value = 'Billy Lindblom'
$ws.Range("AX10").Value = 'Billy Lindblom'

# Row 11
$ws.Range("A11").Value = 111965883
$ws.Range("B11").Value = 55611
$ws.Range("C11").Value = 'Ovaliderad'
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 102612
$ws.Range("F11").Value = 'Järpe'
$ws.Range("G11").Value = 'Tetrastes bonasia'
$ws.Range("H11").Value = '(Linnaeus, 1758)'
$ws.Range("M11").Value = 'lockläte, övriga läten'
$ws.Range("P11").Value = 'Kallhögen 5, Vb'
$ws.Range("Q11").Value = 734846.6442297549
$ws.Range("R11").Value = 7088238.22626837
$ws.Range("S11").Value = 20
$ws.Range("T11").Value = 'Västerbotten'
$ws.Range("U11").Value = 'Vännäs'
$ws.Range("V11").Value = 'Västerbotten'
$ws.Range("W11").Value = 'Vännäs'
$ws.Range("Y11").Value = '''2023-09-08'
$ws.Range("Z11").Value = '17:05'
$ws.Range("AA11").Value = '''2023-09-08'
$ws.Range("AB11").Value = '17:05'
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AW11").Value = 'Billy Lindblom'
$ws.Range("AX11").Value = 'Billy Lindblom'

# Clear M6 (no longer has a value after edit)
$ws.Range("M6").ClearContents()
